$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.441518425941467
$ws.Range("B1").Value = 3.3918616771698
$ws.Range("C1").Value = 4.309640884399414
$ws.Range("D1").Value = 2.128627300262451
$ws.Range("E1").Value = 0.7490930557250977
